$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shared-string text tweaks (report header): volume number + week dates.
#    Edit the longer-length run first is unnecessary here since both target
#    runs keep/alter length in a way that is safe if we patch the later
#    (right-hand) run before the earlier one whose length changes.
# ---------------------------------------------------------------------------

# A8 = "Volume 31   Number  37" -> "...38" (run 4, chars 21-22, same length)
$ws.Range("A8").Characters(21, 2).Text = "38"

# C9 = "Report Covering the Week  9/9/2024  Through  9/15/2024"
#   -> "Report Covering the Week  9/16/2024  Through  9/22/2024"
# Patch the trailing date (same length, 9 chars) first so the earlier date's
# character offset (27) is unaffected when we patch it afterwards.
$ws.Range("C9").Characters(46, 9).Text = "9/22/2024"
$ws.Range("C9").Characters(27, 8).Text = "9/16/2024"

# ---------------------------------------------------------------------------
# Helper "donor" cells used only to copy a cell STYLE (s=14/15/16) across,
# without disturbing their own value - row 14 is untouched by this edit.
# ---------------------------------------------------------------------------
# s="14"  (General / text placeholder style), shared-string "0"     -> C14
# s="14"  (General / text placeholder style), shared-string "***.*" -> E14
# s="15"  (plain integer count style)                               -> I14
# s="16"  (percent-change style)                                    -> K14

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("I14").Copy($ws.Range("C15"))   # switch C15 from text "0" to a number
$ws.Range("C15").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -42.857142857142

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 5
$ws.Range("C14").Copy($ws.Range("C16"))   # switch C16 from a number to text "0"
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -36.842105263157
$ws.Range("I16").Value = 151
$ws.Range("J16").Value = 151
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -13.714285714285
$ws.Range("M16").Value = 51
$ws.Range("N16").Value = -84.654471544715

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 129
$ws.Range("J17").Value = 140
$ws.Range("K17").Value = -7.857142857142
$ws.Range("L17").Value = -3.007518796992
$ws.Range("M17").Value = 86.956521739130
$ws.Range("N17").Value = -40

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -59.090909090909
$ws.Range("I18").Value = 163
$ws.Range("J18").Value = 190
$ws.Range("K18").Value = -14.210526315789
$ws.Range("L18").Value = -10.439560439560
$ws.Range("M18").Value = -10.928961748633
$ws.Range("N18").Value = -92.009803921568

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 47
$ws.Range("E19").Value = -21.276595744680
$ws.Range("F19").Value = 155
$ws.Range("G19").Value = 157
$ws.Range("H19").Value = -1.273885350318
$ws.Range("I19").Value = 1167
$ws.Range("J19").Value = 1266
$ws.Range("K19").Value = -7.819905213270
$ws.Range("L19").Value = -7.011952191235
$ws.Range("M19").Value = 28.665931642778
$ws.Range("N19").Value = -56.127819548872

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 69
$ws.Range("J20").Value = 126
$ws.Range("K20").Value = -45.238095238095
$ws.Range("L20").Value = -50.714285714285
$ws.Range("M20").Value = 13.114754098360
$ws.Range("N20").Value = -97.197400487408

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold row, styles 18/19 - values only change)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = -22.950819672131
$ws.Range("F21").Value = 199
$ws.Range("G21").Value = 226
$ws.Range("H21").Value = -11.946902654867
$ws.Range("I21").Value = 1692
$ws.Range("J21").Value = 1884
$ws.Range("K21").Value = -10.191082802547
$ws.Range("L21").Value = -10.759493670886
$ws.Range("M21").Value = 27.122464312547
$ws.Range("N21").Value = -79.82592106832

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 66.666666666666
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 35
$ws.Range("K22").Value = -2.857142857142
$ws.Range("L22").Value = 3.030303030303
$ws.Range("M22").Value = 54.545454545454

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 75
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 20
$ws.Range("L23").Value = -7.692307692307
$ws.Range("M23").Value = 20

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 58
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = 41.463414634146
$ws.Range("F24").Value = 280
$ws.Range("G24").Value = 228
$ws.Range("H24").Value = 22.807017543859
$ws.Range("I24").Value = 2341
$ws.Range("J24").Value = 2343
$ws.Range("K24").Value = -0.085360648740
$ws.Range("L24").Value = -18.460466736328
$ws.Range("M24").Value = 91.414554374489

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 52
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = 52.941176470588
$ws.Range("F25").Value = 244
$ws.Range("G25").Value = 183
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 2046
$ws.Range("J25").Value = 2065
$ws.Range("K25").Value = -0.920096852300
$ws.Range("L25").Value = -23.370786516853

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -41.666666666666
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 28
$ws.Range("I26").Value = 267
$ws.Range("J26").Value = 236
$ws.Range("K26").Value = 13.135593220339
$ws.Range("L26").Value = -3.956834532374
$ws.Range("M26").Value = 5.952380952380

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("I14").Copy($ws.Range("C27"))   # switch C27 from text "0" to a number
$ws.Range("C27").Value = 3
$ws.Range("F27").Value = 5
$ws.Range("I27").Value = 17
$ws.Range("K27").Value = -10.526315789473
$ws.Range("L27").Value = 30.769230769230

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 12
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 88
$ws.Range("J28").Value = 69
$ws.Range("K28").Value = 27.536231884058
$ws.Range("L28").Value = 41.935483870967

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic. (G/H go from numbers back to the text placeholders)
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("G29"))   # "0"
$ws.Range("E14").Copy($ws.Range("H29"))   # "***.*"

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc. (same placeholder swap as row 29)
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("G30"))   # "0"
$ws.Range("E14").Copy($ws.Range("H30"))   # "***.*"

# ---------------------------------------------------------------------------
# Row 33 - Traffic Fatalities (D/E and G/H go from text placeholders to numbers)
# ---------------------------------------------------------------------------
$ws.Range("I14").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 1
$ws.Range("K14").Copy($ws.Range("E33"))
$ws.Range("E33").Value = -100
$ws.Range("I14").Copy($ws.Range("G33"))
$ws.Range("G33").Value = 1
$ws.Range("K14").Copy($ws.Range("H33"))
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 50
